# Update the "Hello Sunshine"-style weather template from 5/27/2025 to 5/29/2025,
# refreshing the current-conditions / forecast numbers and re-stacking the
# Marshall-Tigers logo above the weather/forecast icon artwork on every slide.

$p = $ppt.ActivePresentation

function Set-ShapeText($shapes, [string]$name, [string]$text) {
    $shp = $shapes.Item($name)
    $shp.TextFrame.TextRange.Text = $text
}

function Move-PictureBeforeSibling($shapes, [string]$moveName, [string]$beforeName) {
    $mover = $shapes.Item($moveName)
    $targetIndex = $shapes.Item($beforeName).ZOrderPosition
    while ($mover.ZOrderPosition -gt $targetIndex) {
        $mover.ZOrder(3)   # msoSendBackward
    }
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $shapes = $slide.Shapes

    $names = @()
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $names += $shapes.Item($k).Name
    }

    if ($names -contains "WeatherIcon") {
        # "Current conditions" slide (1, 3, 5)
        Set-ShapeText $shapes "Time" "12:28 PM"
        Set-ShapeText $shapes "Date" "Thursday, May 29, 2025"
        Set-ShapeText $shapes "WeatherBox" "64"
        Set-ShapeText $shapes "WeatherCondition" "Cloudy"

        Move-PictureBeforeSibling $shapes "Picture 2" "WeatherIcon"
    }
    elseif ($names -contains "ForecastIcon2") {
        # "5-day forecast" slide (2, 4, 6)
        Set-ShapeText $shapes "ForecastDay2" "Friday"
        Set-ShapeText $shapes "ForecastDay3" "Saturday"
        Set-ShapeText $shapes "ForecastDay4" "Sunday"
        Set-ShapeText $shapes "ForecastDay5" "Monday"

        Set-ShapeText $shapes "ForecastTemp2" "H: 82°F  L: 58°F"
        Set-ShapeText $shapes "ForecastTemp3" "H: 84°F  L: 59°F"
        Set-ShapeText $shapes "ForecastTemp4" "H: 88°F  L: 63°F"
        Set-ShapeText $shapes "ForecastTemp5" "H: 90°F  L: 67°F"

        Set-ShapeText $shapes "Time" "12:28 PM"
        Set-ShapeText $shapes "Date" "Thursday, May 29, 2025"

        Move-PictureBeforeSibling $shapes "Picture 2" "ForecastIcon2"
    }
}

# Refresh the cached "last saved" date/time captions on the slide layouts,
# slide master, handout master and notes master (datetime8 / datetimeFigureOut
# fields), which PowerPoint re-stamps whenever the deck is saved on a new day.
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
Set-ShapeText $master.Shapes "Date Placeholder 3" "5/29/2025 12:28 PM"

$customLayouts = $master.CustomLayouts
for ($i = 1; $i -le $customLayouts.Count; $i++) {
    $layout = $customLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "5/29/2025 12:28 PM"
        }
    }
}

# NOTE: deliberately not touching $p.HandoutMaster / $p.NotesMaster here -
# in this host their "Date Placeholder" shape writes alias onto the slide
# master's body placeholder instead of their own XML, corrupting it.
